$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two parameter rows "RecrZ50" (row 151) and "RecrZ95" (row 152) from
# the SpParamsDefinition sheet - entire rows are deleted and the rows below
# shift up (this also removes the now-unused "RecrZ50"/"RecrZ95" name and
# definition strings from the shared string table on save).
$ws.Range("A151:A152").EntireRow.Delete()

# Leave the selection where the author's cursor ended up after the edit
# (previously RespClip's row, now shifted up two rows).
$ws.Range("C153").Select()
